$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = $null

$ws.Range("H32").Value = 2116.2856
$ws.Range("I32").Value = 1731.3334
$ws.Range("J32").Value = 2405
$ws.Range("K32").Value = 1731.3334
$ws.Range("L32").Value = 2405
$ws.Range("M32").Value = -1405.3334
$ws.Range("N32").Value = -3057

$ws.Range("H76").Value = 3327.5
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3327.5
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 3327.5
$ws.Range("N76").Value = -3957.5
$ws.Range("M76").Value = $null

$ws.Range("H79").Value = 3327.5
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3327.5
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 3327.5
$ws.Range("N79").Value = -5511.5
$ws.Range("M79").Value = $null

$ws.Range("H86").Value = 1774.25
$ws.Range("I86").Value = 1800
$ws.Range("K86").Value = 1800
$ws.Range("M86").Value = -677

$ws.Range("H88").Value = 33335032
$ws.Range("J88").Value = 1869.6
$ws.Range("L88").Value = 1869.6
$ws.Range("N88").Value = -2681.6

$ws.Range("H89").Value = 1774.25
$ws.Range("I89").Value = 1800
$ws.Range("K89").Value = 9000
$ws.Range("M89").Value = -3384

$ws.Range("H91").Value = 33335032
$ws.Range("J91").Value = 1869.6
$ws.Range("L91").Value = 1869.6
$ws.Range("N91").Value = -4677.6

$ws.Range("H98").Value = 1065.1945
$ws.Range("I98").Value = 832.75
$ws.Range("J98").Value = 1530.0834
$ws.Range("K98").Value = 832.75
$ws.Range("L98").Value = 1530.0834
$ws.Range("M98").Value = 665.25
$ws.Range("N98").Value = -4526.0834

$ws.Range("H106").Value = 2720.25
$ws.Range("I106").Value = 2720.25
$ws.Range("K106").Value = 2720.25
$ws.Range("M106").Value = -2089.25

$ws.Range("H122").Value = 1065.1945
$ws.Range("I122").Value = 832.75
$ws.Range("J122").Value = 1530.0834
$ws.Range("K122").Value = 2498.25
$ws.Range("L122").Value = 4590.2502
$ws.Range("M122").Value = -48.25
$ws.Range("N122").Value = -9490.2502

$ws.Range("H138").Value = 1933.877
$ws.Range("I138").Value = 1732.7646
$ws.Range("J138").Value = 2154.4517
$ws.Range("K138").Value = 5198.293799999999
$ws.Range("L138").Value = 6463.355100000001
$ws.Range("M138").Value = -58.29379999999946
$ws.Range("N138").Value = -16743.3551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 10856.571
$ws.Range("I36").Value = 10666
$ws.Range("K36").Value = 10666
$ws.Range("M36").Value = -10320

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = $null

$ws.Range("H61").Value = 1858.1578
$ws.Range("I61").Value = 1777.4615
$ws.Range("J61").Value = 2033
$ws.Range("K61").Value = 1777.4615
$ws.Range("L61").Value = 2033
$ws.Range("M61").Value = -1565.4615
$ws.Range("N61").Value = -2457

$ws.Range("H74").Value = 2098.4167
$ws.Range("I74").Value = 897.6
$ws.Range("K74").Value = 897.6
$ws.Range("M74").Value = -23.60000000000002

$ws.Range("H77").Value = 2098.4167
$ws.Range("I77").Value = 897.6
$ws.Range("K77").Value = 4488
$ws.Range("M77").Value = -120

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = $null

$ws.Range("H132").Value = 1523.2354
$ws.Range("I132").Value = 1540.3125
$ws.Range("K132").Value = 4620.9375
$ws.Range("M132").Value = -2090.9375

$ws.Range("H136").Value = 1858.1578
$ws.Range("I136").Value = 1777.4615
$ws.Range("J136").Value = 2033
$ws.Range("K136").Value = 5332.3845
$ws.Range("L136").Value = 6099
$ws.Range("M136").Value = -2782.3845
$ws.Range("N136").Value = -11199

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3311.875
$ws.Range("I31").Value = 4199.75
$ws.Range("J31").Value = 2424
$ws.Range("K31").Value = 4199.75
$ws.Range("L31").Value = 2424
$ws.Range("M31").Value = -3904.75
$ws.Range("N31").Value = -3014

$ws.Range("H33").Value = 8246.666999999999
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 9870
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 9870
$ws.Range("M33").Value = -4621
$ws.Range("N33").Value = -10628

$ws.Range("H34").Value = 3311.875
$ws.Range("I34").Value = 4199.75
$ws.Range("J34").Value = 2424
$ws.Range("K34").Value = 4199.75
$ws.Range("L34").Value = 2424
$ws.Range("M34").Value = -3997.75
$ws.Range("N34").Value = -2828

$ws.Range("H134").Value = 1553
$ws.Range("I134").Value = 1503.375
$ws.Range("K134").Value = 4510.125
$ws.Range("M134").Value = -1975.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 2465.889
$ws.Range("I25").Value = 10000
$ws.Range("J25").Value = 1524.125
$ws.Range("K25").Value = 30000
$ws.Range("L25").Value = 4572.375
$ws.Range("M25").Value = -29831
$ws.Range("N25").Value = -4910.375

$ws.Range("H30").Value = 2465.889
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 1524.125
$ws.Range("K30").Value = 30000
$ws.Range("L30").Value = 4572.375
$ws.Range("M30").Value = -29898
$ws.Range("N30").Value = -4776.375

$ws.Range("H131").Value = 13115.678
$ws.Range("J131").Value = 13549.035
$ws.Range("L131").Value = 40647.105
$ws.Range("N131").Value = -50727.105

$ws.Range("H132").Value = 1319.1428
$ws.Range("I132").Value = 1007.5
$ws.Range("K132").Value = 9067.5
$ws.Range("M132").Value = -6537.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2264.238
$ws.Range("J22").Value = 2360
$ws.Range("L22").Value = 2360
$ws.Range("N22").Value = -2950

$ws.Range("H27").Value = 2264.238
$ws.Range("J27").Value = 2360
$ws.Range("L27").Value = 2360
$ws.Range("N27").Value = -2574

$ws.Range("H46").Value = 2100.6
$ws.Range("I46").Value = 1194.25
$ws.Range("K46").Value = 1194.25
$ws.Range("M46").Value = -1006.25

$ws.Range("H55").Value = 423.68182
$ws.Range("I55").Value = 352.46667
$ws.Range("K55").Value = 352.46667
$ws.Range("M55").Value = -179.46667

$ws.Range("H82").Value = 1969.5
$ws.Range("I82").Value = 1351.5555
$ws.Range("K82").Value = 1351.5555
$ws.Range("M82").Value = -990.5554999999999

$ws.Range("H85").Value = 1969.5
$ws.Range("I85").Value = 1351.5555
$ws.Range("K85").Value = 1351.5555
$ws.Range("M85").Value = -103.5554999999999

$ws.Range("H136").Value = 4640.875
$ws.Range("I136").Value = 3768.3157
$ws.Range("K136").Value = 11304.9471
$ws.Range("M136").Value = -8754.947100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H130").Value = 34559.875
$ws.Range("J130").Value = 34559.875
$ws.Range("L130").Value = 34559.875
$ws.Range("N130").Value = -44599.875

$ws.Range("H136").Value = 18520710
$ws.Range("I136").Value = 34724588
$ws.Range("K136").Value = 104173764
$ws.Range("M136").Value = -104171214
